# Update the cached "datetime1" date-placeholder text from 2024-04-09 to
# 2024-04-23 across the slide master and every slide layout (mirrors using
# Insert > Header & Footer > Apply to All with a fixed date in real
# PowerPoint, which updates every layout/master's cached field text).

$p = $ppt.ActivePresentation

$oldDate = "2024-04-09"
$newDate = "2024-04-23"

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master (날짜 개체 틀 - date placeholder shape).
Update-DateShapes $p.SlideMaster.Shapes

# Every slide layout off of the slide master.
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DateShapes $layouts.Item($L).Shapes
}
